$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Force text storage ("@" number format) for the numeric-looking value columns
# so Excel keeps the literal strings instead of auto-converting to numbers.
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:F4").NumberFormat = "@"

# Row J_0_L0_v
$ws.Range("A2").Value = "2.9999999999999996 - 2x_1 + y_1 - y_2"
$ws.Range("B2").Value = "-0.49999999999999956"
$ws.Range("D2").Value = "0.96"
$ws.Range("E2").Value = "5.8"
$ws.Range("F2").Value = "5.6000000000000005"

# Row J_0_LP_v
$ws.Range("A3").Value = "-0.9499999999999997 + x_1 - 3x_2 + y_2"
$ws.Range("B3").Value = "-1.0500000000000003"
$ws.Range("D3").Value = "0.9"
$ws.Range("E3").Value = "8.4"
$ws.Range("F3").Value = "3.5"

# Row J_Ne_L0_v
$ws.Range("A4").Value = "-4.79 + x_1 + x_2"
$ws.Range("B4").Value = "1.7999999999999998"
$ws.Range("D4").Value = "0.28"
$ws.Range("E4").Value = "6.4"
$ws.Range("F4").Value = "4.4"

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2:D2").NumberFormat = "@"
$ws.Range("A2").Value = "2.55"
$ws.Range("B2").Value = "1.25"
$ws.Range("C2").Value = "4.25"
$ws.Range("D2").Value = "2.15"

# --- Sheet: Vector_bf ---
# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"); Worksheets.Item(name) resolves
# case-insensitively (like real Excel) so both names would hit the same
# sheet. Use the unambiguous 1-based tab index instead (5 = Vector_bf).
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "3.04"
$ws.Range("A3").Value = "-0.9400000000000001"

# --- Sheet: Vector_BF --- (tab index 6 = Vector_BF)
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2").Value = "-1.200000000000001"
$ws.Range("A3").Value = "17.800000000000004"
$ws.Range("A4").Value = "-6.3"
$ws.Range("A5").Value = "-2.6000000000000005"
